$d = $word.ActiveDocument

# ============================================================
# Strategist section rewrite (V15 - "Strategist gebruiken is geen actie")
# ============================================================

# --- Paragraph 1 --------------------------------------------------
# "After revealing, the Strategist can reveal any opposing piece that it
#  can reach, meaning it forces the owner of that piece to reveal the
#  piece's character."
# becomes
# "At the beginning of a turn, the Strategist can reveal itself and
#  inspect any opposing piece that it can reach; the owner of that piece
#  must reveal the piece's character."

$d.Content.Find.Execute(
    "After revealing, the Strategist can reveal",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "At the beginning of a turn, the Strategist can reveal",
    2) | Out-Null

$d.Content.Find.Execute(
    "reveal any opposing piece that it can reach, meaning it forces ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "reveal itself and inspect any opposing piece that it can reach; ",
    2) | Out-Null

$d.Content.Find.Execute(
    "the owner of that piece to reveal the piece",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "the owner of that piece must reveal the piece",
    2) | Out-Null

# --- Paragraph 2 --------------------------------------------------
# "This is done instead of moving or capturing."
# becomes
# "This is done before the player moves or captures pieces. The
#  Strategist can only inspect one piece per turn."

$d.Content.Find.Execute(
    "This is done instead of moving or capturing.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This is done before the player moves or captures pieces. The Strategist can only inspect one piece per turn.",
    2) | Out-Null

# --- Paragraph 3 --------------------------------------------------
# "...the Strategist can reveal a piece in the Garden. As an Elephant,
#  the Strategist cannot reveal pieces more than two spaces away."
# becomes
# "...the Strategist can inspect a piece in the Garden. As an Elephant,
#  the Strategist cannot inspect pieces more than two spaces away."

$d.Content.Find.Execute(
    "the Strategist can reveal a piece in the Garden",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "the Strategist can inspect a piece in the Garden",
    2) | Out-Null

$d.Content.Find.Execute(
    "the Strategist cannot reveal pieces more than two spaces away",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "the Strategist cannot inspect pieces more than two spaces away",
    2) | Out-Null

# ============================================================
# Move the hidden "_GoBack" (last-edit-location) bookmark from the end
# of the Elephant paragraph to right after "cannot inspect" in the
# Strategist's final paragraph, matching where the real edit last
# touched the document.
# ============================================================

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$target = $d.Content
$target.Find.Execute(
    "As an Elephant, the Strategist cannot inspect",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "", 0) | Out-Null
$target.Collapse(0)  # wdCollapseEnd
$d.Bookmarks.Add("_GoBack", $target) | Out-Null
